$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Lösung für die Performance-Optimierung" -> "Caching Solution"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lösung für die Performance-Optimierung", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Caching Solution", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Problembeschreibun" + bookmark("_GoBack") + "g" -> single run "Problembeschreibung"
#    (searching across the full word consumes the in-between bookmark too)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Problembeschreibung", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Problembeschreibung", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Merge "Anhand der Risiken mussten wir unsere Möglichkeit selektieren" + ", "
#    into one run/text chunk (the trailing "dahingehend ... untersucht." keeps reading the same).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Anhand der Risiken mussten wir unsere Möglichkeit selektieren, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Anhand der Risiken mussten wir unsere Möglichkeit selektieren, ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Auf dem Server cachen" list item: bold + dark red (C00000), plus append
#    " (kurz Beschreibung bitte)".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Auf dem Server cachen*") {
        $insertPos = $p.Range.End - 1
        $ins = $d.Range($insertPos, $insertPos)
        $ins.InsertAfter(" (kurz Beschreibung bitte)")

        # Bold + red colour across the whole paragraph (existing + new text)
        $p.Range.Font.Bold = 1
        $p.Range.Font.Color = 192   # wdColor value for RGB C00000
        $p.Range.Font.Underline = 11 # wdUnderlineWavy -> <w:u w:val="wave"/>
        break
    }
}

# ---------------------------------------------------------------------------
# 5. Move the "_GoBack" bookmark so it now sits right before the
#    "Effizientere Software schreiben" paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Effizientere Software schreiben*") {
        $d.Bookmarks.Add("_GoBack", $d.Range($p.Range.Start, $p.Range.Start)) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 6. "„Outsourcing“" list item: bold + dark red (C00000), plus append
#    " (kurz Beschreibung bitte)".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Outsourcing*") {
        $insertPos = $p.Range.End - 1
        $ins = $d.Range($insertPos, $insertPos)
        $ins.InsertAfter(" (kurz Beschreibung bitte)")

        $p.Range.Font.Bold = 1
        $p.Range.Font.Color = 192
        $p.Range.Font.Underline = 11
        break
    }
}
